$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (instead of numeric auto-detection) for price cells
# whose new values would otherwise be parsed as plain numbers (losing the
# original text formatting, e.g. trailing zeros such as "21.10" -> 21.1).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "37.282.60"
$ws.Range("E2").Value = "  +2.07%  "

# Row 3
$ws.Range("D3").Value = "2.061.53"
$ws.Range("E3").Value = "  +3.36%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "234.35"
$ws.Range("E5").Value = "  -1.02%  "

# Row 6
$ws.Range("E6").Value = "  +2.43%  "

# Row 7
$ws.Range("E7").Value = "  +5.79%  "

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("E9").Value = "  +2.67%  "

# Row 10
$ws.Range("D10").Value = "58.72"
$ws.Range("E10").Value = "  +1.57%  "

# Row 11
$ws.Range("E11").Value = "  +1.64%  "

# Row 12
$ws.Range("E12").Value = "  +2.75%  "

# Row 13
$ws.Range("D13").Value = "2.364.73"
$ws.Range("E13").Value = "  +3.42%  "

# Row 14
$ws.Range("E14").Value = "  +2.88%  "

# Row 15
$ws.Range("D15").Value = "21.10"
$ws.Range("E15").Value = "  +3.12%  "

# Row 16
$ws.Range("D16").Value = "0.774"
$ws.Range("E16").Value = "  +2.23%  "

# Row 17
$ws.Range("E17").Value = "  +1.77%  "

# Row 18
$ws.Range("D18").Value = "2.062.32"
$ws.Range("E18").Value = "  +3.04%  "

# Row 19
$ws.Range("D19").Value = "37.375.33"
$ws.Range("E19").Value = "  +2.40%  "

# Row 20
$ws.Range("D20").Value = "6.15"
$ws.Range("E20").Value = "  +16.85%  "

# Row 21
$ws.Range("D21").Value = "69.45"
$ws.Range("E21").Value = "  +2.36%  "

# Row 22
$ws.Range("E22").Value = "  +0.55%  "

# Row 23
$ws.Range("D23").Value = "226.93"
$ws.Range("E23").Value = "  +2.25%  "

# Row 24
$ws.Range("E24").Value = "  +0.09%  "

# Row 25
$ws.Range("E25").Value = "  +1.27%  "

# Row 26
$ws.Range("E26").Value = "  +0.95%  "

# Row 27
$ws.Range("D27").Value = "165.72"
$ws.Range("E27").Value = "  +2.19%  "

# Row 28
$ws.Range("D28").Value = "1.50"
$ws.Range("E28").Value = "  +11.57%  "

# Row 29
$ws.Range("E29").Value = "  +1.99%  "

# Row 30
$ws.Range("D30").Value = "19.15"
$ws.Range("E30").Value = "  +1.48%  "

# Row 31
$ws.Range("E31").Value = "  -1.58%  "

# Row 32
$ws.Range("E32").Value = "  +1.80%  "

# Row 33
$ws.Range("D33").Value = "4.52"
$ws.Range("E33").Value = "  +3.25%  "

# Row 34
$ws.Range("D34").Value = "0.0621"
$ws.Range("E34").Value = "  +2.33%  "

# Row 35
$ws.Range("E35").Value = "  +8.68%  "

# Row 36
$ws.Range("E36").Value = "  +6.73%  "

# Row 37
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "3.38"
$ws.Range("E37").Value = "  -0.96%  "

# Row 38
$ws.Range("E38").Value = "  -0.05%  "

# Row 39
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").Value = "1.80"
$ws.Range("E39").Value = "  +1.76%  "

# Row 40
$ws.Range("E40").Value = "  +4.55%  "

# Row 41
$ws.Range("D41").Value = "0.0982"
$ws.Range("E41").Value = "  +3.81%  "

# Row 42
$ws.Range("E42").Value = "  -1.31%  "

# Row 43
$ws.Range("D43").Value = "4.36"
$ws.Range("E43").Value = "  +23.14%  "

# Row 44
$ws.Range("D44").Value = "1.456.30"
$ws.Range("E44").Value = "  +0.28%  "

# Row 45
$ws.Range("D45").Value = "95.39"
$ws.Range("E45").Value = "  +6.88%  "

# Row 46
$ws.Range("E46").Value = "  +3.66%  "

# Row 47
$ws.Range("E47").Value = "  +4.09%  "

# Row 48
$ws.Range("D48").Value = "15.81"
$ws.Range("E48").Value = "  +3.92%  "

# Row 49
$ws.Range("E49").Value = "  +3.27%  "

# Row 50
$ws.Range("D50").Value = "7.26"
$ws.Range("E50").Value = "  +5.66%  "

# Row 51
$ws.Range("E51").Value = "  +2.01%  "
